$wb = $excel.ActiveWorkbook

$oldName   = "572775a6-2d63-486d-8196-f83075dc5894.md"
$newName   = "875889ae-a1ac-4d6a-8e32-44fa8d7e098e.md"
$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"
$epoch     = "0001-01-01 00:00:00"

function Update-HeadHyperlink($ws) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.TextToDisplay -eq $oldName) {
            $h.TextToDisplay = $newName
        }
    }
}

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = $newName
$ws.Range("B2").Value = $newStatus
$ws.Range("C2").Value = $newStatus
Update-HeadHyperlink $ws

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = $newName
$ws.Range("B2").Value = $newStatus
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = $epoch
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Ignored"
foreach ($h in $ws.Hyperlinks) {
    if ($h.TextToDisplay -eq $oldName) {
        $h.TextToDisplay = $newName
    } elseif ($h.Range.Address -eq $ws.Range("C2").Address) {
        $h.Delete()
    }
}

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = $newName
$ws.Range("B2").Value = $newStatus
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = $epoch
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Ignored"
foreach ($h in $ws.Hyperlinks) {
    if ($h.TextToDisplay -eq $oldName) {
        $h.TextToDisplay = $newName
    } elseif ($h.Range.Address -eq $ws.Range("C2").Address) {
        $h.Delete()
    }
}
